# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.311.54'
$ws.Range('D3').Value = '1.872.66'
$ws.Range('D4').Value = "'" + '1.001'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('E5').Value = '  -0.81%  '
$ws.Range('D6').Value = "'" + '241.86'
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('D7').Value = "'" + '1.001'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = "'" + '0.07802'
$ws.Range('E8').Value = '  +1.18%  '
$ws.Range('D9').Value = "'" + '0.3105'
$ws.Range('E9').Value = '  -0.33%  '
$ws.Range('D10').Value = "'" + '25.10'
$ws.Range('E10').Value = '  -0.22%  '
$ws.Range('D11').Value = "'" + '0.08376'
$ws.Range('D12').Value = '1.886.39'
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D13').Value = "'" + '5.235'
$ws.Range('E13').Value = '  -0.07%  '
$ws.Range('D14').Value = "'" + '0.7166'
$ws.Range('E14').Value = '  +0.32%  '
$ws.Range('E15').Value = '  -0.48%  '
$ws.Range('D16').Value = "'" + '0.000008388'
$ws.Range('E16').Value = '  +0.64%  '
$ws.Range('D17').Value = "'" + '6.134'
$ws.Range('E17').Value = '  +2.85%  '
$ws.Range('D18').Value = '29.318.06'
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('D19').Value = "'" + '240.31'
$ws.Range('E19').Value = '  -1.21%  '
$ws.Range('D20').Value = '2.126.97'
$ws.Range('E20').Value = '  -0.73%  '
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('D22').Value = "'" + '1.001'
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('D23').Value = "'" + '7.738'
$ws.Range('E23').Value = '  -1.80%  '
$ws.Range('D24').Value = "'" + '1.001'
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').Value = "'" + '0.1595'
$ws.Range('E25').Value = '  -1.46%  '
$ws.Range('D26').Value = "'" + '162.60'
$ws.Range('E26').Value = '  -0.74%  '
$ws.Range('D27').Value = "'" + '9.031'
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').Value = "'" + '18.47'
$ws.Range('E28').Value = '  -0.45%  '
$ws.Range('D29').Value = "'" + '1.505'
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('D30').Value = "'" + '4.416'
$ws.Range('E30').Value = '  +0.28%  '
$ws.Range('D31').Value = "'" + '4.352'
$ws.Range('E31').Value = '  +0.55%  '
$ws.Range('D32').Value = "'" + '1.231'
$ws.Range('E32').Value = '  -4.64%  '
$ws.Range('D33').Value = "'" + '0.05355'
$ws.Range('E33').Value = '  +1.92%  '
$ws.Range('E34').Value = '  +0.98%  '
$ws.Range('D35').Value = "'" + '0.7486'
$ws.Range('E35').Value = '  -1.46%  '
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('D37').Value = "'" + '2.685'
$ws.Range('E37').Value = '  +0.25%  '
$ws.Range('D38').Value = "'" + '0.01880'
$ws.Range('E38').Value = '  +0.86%  '
$ws.Range('D39').Value = '1.242.28'
$ws.Range('E39').Value = '  +6.85%  '
$ws.Range('D40').Value = "'" + '2.732'
$ws.Range('E40').Value = '  +0.40%  '
$ws.Range('D41').Value = "'" + '6.531'
$ws.Range('E41').Value = '  +2.82%  '
$ws.Range('D42').Value = "'" + '0.8921'
$ws.Range('E42').Value = '  +0.29%  '
$ws.Range('D43').Value = "'" + '109.91'
$ws.Range('E43').Value = '  +4.94%  '
$ws.Range('D44').Value = "'" + '72.26'
$ws.Range('E44').Value = '  -1.75%  '
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '2.019.58'
$ws.Range('E46').Value = '  -0.49%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = "'" + '0.5201'
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = "'" + '1.794'
$ws.Range('E48').Value = '  -0.29%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = "'" + '9.457'
$ws.Range('E49').Value = '  +0.65%  '
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D50').Value = "'" + '0.4336'
$ws.Range('E50').Value = '  +0.69%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').Value = "'" + '7.095'
$ws.Range('E51').Value = '  +0.70%  '
